# The commit swaps the deck's theme so that the "Office Theme" colour
# palette (previously the inactive/second theme) becomes the palette the
# presentation actually uses (ppt/theme/theme1.xml, the theme wired to the
# slide master / presentation root), while the old "Integral" palette moves
# off to the side.
#
# The PowerPoint object model doesn't give a way to overwrite a whole
# <a:theme> part at once (no file-system access inside this host), so the
# supported way to edit a theme is per theme colour, via
# ThemeColorScheme.Colors(i).RGB - exactly what we do below, writing every
# slot of the 12-colour scheme (dk1, lt1, dk2, lt2, accent1-6, hlink,
# folHlink) with the Office Theme's values.

$p = $ppt.ActivePresentation
$theme = $p.SlideMaster.Theme
$tcs = $theme.ThemeColorScheme

# Office Theme colour scheme values, in clrScheme slot order
# (dk1, lt1, dk2, lt2, accent1..6, hlink, folHlink).
$officeColors = @(
    0x000000,  # dk1
    0xFFFFFF,  # lt1
    0x44546A,  # dk2
    0xE7E6E6,  # lt2
    0x5B9BD5,  # accent1
    0xED7D31,  # accent2
    0xA5A5A5,  # accent3
    0xFFC000,  # accent4
    0x4472C4,  # accent5
    0x70AD47,  # accent6
    0x0563C1,  # hlink
    0x954F72   # folHlink
)

for ($i = 1; $i -le $tcs.Count; $i++) {
    $hex = $officeColors[$i - 1]
    $r = [math]::Floor($hex / 0x10000) -band 0xFF
    $g = [math]::Floor($hex / 0x100) -band 0xFF
    $b = $hex -band 0xFF
    # PowerPoint's RGB/ColorFormat.RGB integer is r + g*256 + b*65536.
    $tcs.Item($i).RGB = $r + ($g * 256) + ($b * 65536)
}
